$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 4 (mirrors the pattern of rows 2/3, new test case "3") ---
$ws.Range("A4").Value = "'3"
$ws.Range("B4").Value = "'95400152"
$ws.Range("B4").Font.Size = 12
$ws.Range("C4").Value = "'1"
$ws.Range("D4").Value = "'sandrita69"
$ws.Range("D4").Font.Size = 12
$ws.Range("E4").Value = "'1234"
$ws.Range("E4").Font.Size = 12
$ws.Range("F4").Value = "'4321"
$ws.Range("F4").Font.Size = 12
$ws.Range("G4").Value = "Acierto"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Borders.LineStyle = 1
$ws.Range("H4").Value = "'001"
$ws.Range("I4").Value = "'0370"
$ws.Range("J4").Value = "'NO ERROR"
$ws.Range("K4").Value = "'bolp"
$ws.Range("L4").Value = "'ACTIVO"

$ws.Rows.Item(4).RowHeight = 15.75

# --- Update G2 (orientacion) ---
$ws.Range("G2").Value = "Alterno"

# --- Data validation dropdown list for G4:J4 ---
$ws.Range("G4:J4").Validation.Add(3, 1, 1, "[1]Listas!#REF!")

# --- Selection moves to G4 ---
$ws.Range("G4").Select()
